$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the header label for column B ("value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# 2. Extend the date-formatted style (copied from an existing date cell) down
#    over the newly-added rows before writing values into them.
$ws.Range("A2").Copy()
$ws.Range("A54:A84").PasteSpecial(-4122)

# 3. Write the full refreshed date/value series (rows 2-84).
$dates = @(38398, 38487, 38579, 38671, 38763, 38852, 38944, 39036, 39128, 39217, 39309, 39401, 39493, 39583, 39675, 39767, 39859, 39948, 40040, 40132, 40224, 40313, 40405, 40497, 40589, 40678, 40770, 40862, 40954, 41044, 41136, 41228, 41320, 41409, 41501, 41593, 41685, 41774, 41866, 41958, 42050, 42139, 42231, 42323, 42415, 42505, 42597, 42689, 42781, 42870, 42962, 43054, 43146, 43235, 43327, 43419, 43511, 43600, 43692, 43784, 43876, 43966, 44058, 44150, 44242, 44331, 44423, 44515, 44607, 44696, 44788, 44880, 44972, 45061, 45153, 45245, 45337, 45427, 45519, 45611, 45703, 45792, 45884)
$values = @(0.4001302730732021, 1.534309226294653, -0.2158762754026498, -0.1573369388209471, 0.5441785023706558, -0.1273572285275435, 0.510068525034896, -0.009760955203091726, 1.786050651751793, -0.03834288659695062, -0.0575559711994913, -0.6200525020039009, 1.257433230729447, 0.8837904892317567, 0.810806491930748, -0.07476001263452758, 0.2524978494830066, 0.5778148852415939, 0.1185568564730346, -0.1977656654399595, 1.071871022829441, -1.091011900795806, 1.120967691003898, -0.1775928823643795, 1.28981182300268, 0.5923450763659872, 0.5531759638372762, 0.5788603179058356, 0.1947850960503388, -0.2481858862331165, 0.3732050716642448, 0.141299961337424, -0.1411005862636046, -0.2128461555332564, 0.4532479246724535, -0.3446087745608111, 0.4255979180752121, 0.4461687925667093, 0.6373066379050414, 0.3261422475203943, 0.6788370390783598, 0.6601374471387373, 1.255382587579845, 0.8728685839363095, 0.4694885089849095, 1.152137745180852, 0.9596379771730028, 0.2682953781150843, 0.4191917022489378, 0.1682050168937224, -0.03534872415686152, 0.5002605909365485, -0.5240674734835977, 0.5854015665873362, 0.2, 1.628071843823122, -0.3, 0.5022917647287812, 0.754883892913071, 0.3494637214130449, 0.2, 1.450185044412038, 0.3499999990000049, 0.700000000000017, 0.700000000000017, 0.7999999999999972, -0.7999999999999972, 0.4999999999999858, 0.00000000000002842170943040401, 0.5000000000000142, 0.700000000000017, -0.9999999999999858, -0.5, 2.799999999999997, 0.0, 0.2000000000000028, -0.09999999999999432, 0.2000000000000028, 0.09999999999999432, 0.09999999999999432, 0.0, 0.4999999999999858, 0.4000000000000199)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $values[$i]
}
